$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content edits: row 12 (R8) ---
$ws.Range("B12").Value = "Requiere licencia de Power BI."
$ws.Range("C12").Value = "Revisión."
$ws.Range("D12").Value = "Cumple."
$ws.Range("E12").Value = "Licencia Power BI, cumple."

# --- Content edits: row 14 (R10) ---
$ws.Range("C14").Value = "Revisión."
$ws.Range("D14").Value = "Cumple."
$ws.Range("E14").Value = "Filtros disponibles para segmentar por clientes y sectores principalmente, cumple."

# --- View state: zoom + selection ---
$win = $excel.ActiveWindow
$win.Zoom = 80
[void]$ws.Range("E15").Select()
